$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric-value-only updates (column C) ---
$ws.Range("C3").Value = 1051.0
$ws.Range("C4").Value = 735.0
$ws.Range("C5").Value = 469.0
$ws.Range("C10").Value = 370.0
$ws.Range("C12").Value = 307.0
$ws.Range("C13").Value = 232.0
$ws.Range("C15").Value = 165.0
$ws.Range("C16").Value = 144.0
$ws.Range("C19").Value = 126.0
$ws.Range("C32").Value = 44.0

# --- Rows 8-9: swap entities (with updated counts) ---
$ws.Range("A8").Value = "Regione Lombardia"
$ws.Range("B8").Formula = "'80050050154"
$ws.Range("C8").Value = 374.0
$ws.Range("A9").Value = "Lepida Spa"
$ws.Range("B9").Formula = "'02770891204"
$ws.Range("C9").Value = 372.0

# --- Rows 24-26: rotate three entities (with updated counts) ---
$ws.Range("A24").Value = "Regione Piemonte"
$ws.Range("B24").Formula = "'80087670016"
$ws.Range("C24").Value = 79.0
$ws.Range("A25").Value = "Next Step Solution"
$ws.Range("B25").Formula = "'02554480349"
$ws.Range("C25").Value = 74.0
$ws.Range("A26").Value = "Regione Autonoma Friuli-Venezia Giulia"
$ws.Range("B26").Formula = "'80014930327"
$ws.Range("C26").Value = 73.0

# --- Rows 38-39: swap entities (with updated counts) ---
$ws.Range("A38").Value = "PMPay s.r.l."
$ws.Range("B38").Formula = "'08747230962"
$ws.Range("C38").Value = 33.0
$ws.Range("A39").Value = "ROMA CAPITALE"
$ws.Range("B39").Formula = "'02438750586"
$ws.Range("C39").Value = 32.0

# --- Rows 46-47: swap entities (counts unchanged) ---
$ws.Range("A46").Value = "ANDREANI TRIBUTI srl"
$ws.Range("B46").Formula = "'01412920439"
$ws.Range("A47").Value = "Comune di Palermo"
$ws.Range("B47").Formula = "'80016350821"

# --- Rows 55 & 57: swap entities (counts unchanged) ---
$ws.Range("A55").Value = "Comune di Catania"
$ws.Range("B55").Formula = "'00137020871"
$ws.Range("A57").Value = "Numera Sistemi e Informatica SpA"
$ws.Range("B57").Formula = "'01265230902"

# --- Rows 58-72: names/codes shift by one position (counts unchanged) ---
$ws.Range("A58").Value = "Linea Comune Spa"
$ws.Range("B58").Formula = "'05591950489"
$ws.Range("A59").Value = "ISWEB S.p.A."
$ws.Range("B59").Formula = "'01722270665"
$ws.Range("A60").Value = "I.C.A. - Imposte Comunali Affini – s.r.l."
$ws.Range("B60").Formula = "'02478610583"
$ws.Range("A61").Value = "ICCREA Banca SpA"
$ws.Range("B61").Formula = "'04774801007"
$ws.Range("A62").Value = "Engineering Ingegneria Informatica SpA"
$ws.Range("B62").Formula = "'00967720285"
$ws.Range("A63").Value = "Ministero dello Sviluppo Economico"
$ws.Range("B63").Formula = "'80230390587"
$ws.Range("A64").Value = "Softline srl"
$ws.Range("B64").Formula = "'12299030150"
$ws.Range("A65").Value = "CityPoste Payment Digital S.r.l."
$ws.Range("B65").Formula = "'02003750672"
$ws.Range("A66").Value = "Agenzia Italiana del Farmaco - AIFA"
$ws.Range("B66").Formula = "'97345810580"
$ws.Range("A67").Value = "Società Almaviva S.p.A."
$ws.Range("B67").Formula = "'08450891000"
$ws.Range("A68").Value = "Banco BPM Società per Azioni"
$ws.Range("B68").Formula = "'09722490969"
$ws.Range("A69").Value = "ARGO SOFTWARE SRL"
$ws.Range("B69").Formula = "'00838520880"
$ws.Range("A70").Value = "MegASP S.r.l."
$ws.Range("B70").Formula = "'09898030151"
$ws.Range("A71").Value = "ARCA Servizi s.r.l"
$ws.Range("B71").Formula = "'09106071005"
$ws.Range("A72").Value = "San Marco SPA"
$ws.Range("B72").Formula = "'04142440728"
